# Add a header row ("Code" / "State") above the existing state-code table.
# All existing data (codes in column A, state names in column B) shifts
# down by one row automatically when the new row is inserted.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at the top; existing rows 1..62 become 2..63.
$ws.Rows("1:1").Insert()

# Populate the new header row.
$ws.Range("A1").Value = "Code"
$ws.Range("B1").Value = "State"

# Match the saved selection state (active cell on A2).
$ws.Range("A2").Select()
